# Update cryptocurrency price and 1h-volume figures in the cryptos worksheet
# (refresh snapshot committed by the scheduled GitHub Actions scraper run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.291.16"
$ws.Range("E2").Value = "  -4.02%  "
$ws.Range("D3").Value = "1.664.01"
$ws.Range("E3").Value = "  -2.71%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.59"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5158"
$ws.Range("E6").Value = "  -3.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06433"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2568"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.86"
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07658"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.669.11"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").Value = "1.895.28"
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.299"
$ws.Range("E14").Value = "  -5.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5537"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").Value = "0.0₅8028"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.43"
$ws.Range("E17").Value = "  -5.14%  "
$ws.Range("D18").Value = "26.344.31"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.03"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.388"
$ws.Range("E21").Value = "  -6.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.10"
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.888"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.16"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.755"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1159"
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.958"
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.75"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05252"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.261"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.369"
$ws.Range("E32").Value = "  -4.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.222"
$ws.Range("E33").Value = "  -6.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.563"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.755"
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9213"
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5724"
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("D39").Value = "1.152.95"
$ws.Range("E39").Value = "  +10.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01594"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.009"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8407"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.648"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.85"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "1.804.86"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("E46").Value = "  -7.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4505"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.01"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.909"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05100"
$ws.Range("E51").Value = "  -2.71%  "
